# Update "想去人数" (want-to-go count) figures in column F across sheets,
# reflecting refreshed counts from the upstream data source.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1417
$ws1.Range("F5").Value = 6509
$ws1.Range("F6").Value = 508
$ws1.Range("F8").Value = 22
$ws1.Range("F10").Value = 6732
$ws1.Range("F12").Value = 1370
$ws1.Range("F13").Value = 789
$ws1.Range("F23").Value = 1024
$ws1.Range("F24").Value = 331
$ws1.Range("F39").Value = 304

# Sheet 2: 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F22").Value = 187
$ws2.Range("F27").Value = 610
$ws2.Range("F31").Value = 719
$ws2.Range("F41").Value = 56

# Sheet 3: 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F6").Value = 583
$ws3.Range("F8").Value = 1185

# Sheet 4: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 583
$ws4.Range("F10").Value = 6509
$ws4.Range("F11").Value = 508
$ws4.Range("F13").Value = 22
$ws4.Range("F15").Value = 6732
$ws4.Range("F18").Value = 1370
$ws4.Range("F24").Value = 1185
$ws4.Range("F27").Value = 187
$ws4.Range("F31").Value = 610
$ws4.Range("F46").Value = 304
$ws4.Range("F50").Value = 56
